# Fruta / hortaliza, semanal
#
# The data rows (2-17) on the sheet got reshuffled: every row's full set of
# field values (date, quality, volume, prices, unit, origin, etc.) moved to
# a different row in the table; row 4 stayed put. We snapshot the whole
# A2:T17 block first (so the permutation can be applied safely, independent
# of write order) and then write each destination row from its mapped
# source row in that snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17
$firstCol = 1
$lastCol = 20

$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$snapshot = $srcRange.Value2()

# destRow (1-based within the block, i.e. Excel row - 1) -> sourceRow (same indexing)
$rowMap = @{
  1  = 14   # Excel row 2  <- old row 15
  2  = 9    # Excel row 3  <- old row 10
  3  = 3    # Excel row 4  <- old row 4 (unchanged)
  4  = 8    # Excel row 5  <- old row 9
  5  = 16   # Excel row 6  <- old row 17
  6  = 15   # Excel row 7  <- old row 16
  7  = 5    # Excel row 8  <- old row 6
  8  = 1    # Excel row 9  <- old row 2
  9  = 13   # Excel row 10 <- old row 14
  10 = 6    # Excel row 11 <- old row 7
  11 = 10   # Excel row 12 <- old row 11
  12 = 11   # Excel row 13 <- old row 12
  13 = 12   # Excel row 14 <- old row 13
  14 = 7    # Excel row 15 <- old row 8
  15 = 4    # Excel row 16 <- old row 5
  16 = 2    # Excel row 17 <- old row 3
}

$numRows = $lastRow - $firstRow + 1
$numCols = $lastCol - $firstCol + 1

$newValues = New-Object 'object[,]' $numRows, $numCols

for ($destIdx = 1; $destIdx -le $numRows; $destIdx++) {
  $srcIdx = $rowMap[$destIdx]
  for ($col = 1; $col -le $numCols; $col++) {
    $newValues[$destIdx - 1, $col - 1] = $snapshot[$srcIdx, $col]
  }
}

$srcRange.Value = $newValues
